# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

$wsYDS   = $wb.Worksheets.Item("YDS")
$wsOFF   = $wb.Worksheets.Item("OFF")
$wsDEF   = $wb.Worksheets.Item("DEF")
$wsST    = $wb.Worksheets.Item("ST")
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsPEN   = $wb.Worksheets.Item("PEN")

# ---- YDS sheet: append this week's per-drive(ish) yardage logs ----
$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value() + " 43 3 -1 4 1 -2 2 2 4 -2 7 -2 5 5 11 -1 2 8 18 4 3 -1 1 3 3 -1 3 2 2 1 -3 3"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value() + " 9 1 6 23 10 7 3 12 14 7 5 7 7 5 5 37 8 20 39 14"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value() + " 6 5 4 2 2 10 57 0 2 1 4 1 0 5 1 7 0 11 10 2"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value() + " 3 10 5 10 5 8 13 7 11 4 3 11 10 24 3 33 4 5 7 16 7 4 7 14 2 14 8"

# ---- OFF sheet: updated season totals (Home row 2, Road row 3) ----
$wsOFF.Range("C2").Value = 415
$wsOFF.Range("D2").Value = 28
$wsOFF.Range("E2").Value = 17
$wsOFF.Range("F2").Value = 152
$wsOFF.Range("G2").Value = 122
$wsOFF.Range("H2").Value = 7
$wsOFF.Range("J2").Value = 73
$wsOFF.Range("N2").Value = 25
$wsOFF.Range("O2").Value = 51

$wsOFF.Range("B3").Value = 20
$wsOFF.Range("C3").Value = 340
$wsOFF.Range("E3").Value = 67
$wsOFF.Range("F3").Value = 186
$wsOFF.Range("H3").Value = 48
$wsOFF.Range("I3").Value = 113
$wsOFF.Range("J3").Value = 92
$wsOFF.Range("L3").Value = 523
$wsOFF.Range("M3").Value = 347
$wsOFF.Range("Q3").Value = 1034

# ---- DEF sheet: updated season totals (Home row 2, Road row 3) ----
$wsDEF.Range("C2").Value = 363
$wsDEF.Range("D2").Value = 21
$wsDEF.Range("F2").Value = 110
$wsDEF.Range("G2").Value = 107
$wsDEF.Range("J2").Value = 61
$wsDEF.Range("O2").Value = 32
$wsDEF.Range("P2").Value = 15

$wsDEF.Range("B3").Value = 15
$wsDEF.Range("C3").Value = 389
$wsDEF.Range("E3").Value = 56
$wsDEF.Range("F3").Value = 220
$wsDEF.Range("G3").Value = 72
$wsDEF.Range("I3").Value = 124
$wsDEF.Range("J3").Value = 104
$wsDEF.Range("L3").Value = 568
$wsDEF.Range("M3").Value = 374
$wsDEF.Range("Q3").Value = 954

# ---- ST sheet: updated season totals + appended per-kick logs ----
$wsST.Range("B2").Value = 185
$wsST.Range("D2").Value = 101
$wsST.Range("F2").Value = 129
$wsST.Range("G2").Value = 124
$wsST.Range("J2").Value = 42
$wsST.Range("K2").Value = 41
$wsST.Range("L2").Value = 33
$wsST.Range("M2").Value = 24
$wsST.Range("N2").Value = 10

$wsST.Range("B3").Value = 88

$wsST.Range("B4").Value = $wsST.Range("B4").Value() + " 66 68 49 68"
$wsST.Range("B5").Value = $wsST.Range("B5").Value() + " 21 13 16 31"
$wsST.Range("B6").Value = $wsST.Range("B6").Value() + " 16"
$wsST.Range("D3").Value = $wsST.Range("D3").Value() + " 39 44 32"
$wsST.Range("D4").Value = $wsST.Range("D4").Value() + " 0 48 0"
$wsST.Range("D5").Value = $wsST.Range("D5").Value() + " 0"

# ---- TURNS sheet: updated season totals (Road row 3) ----
$wsTURNS.Range("D3").Value = 11
$wsTURNS.Range("E3").Value = 16

# ---- PEN sheet: updated season totals ----
$wsPEN.Range("B2").Value = 26
$wsPEN.Range("B3").Value = 40
